$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed crypto price / volume(1h) figures scraped this run.
# A handful of "Price" cells are bare decimal numbers (e.g. "1.007") which Excel
# would otherwise auto-convert to a numeric value; prefix those with a literal
# leading apostrophe (the same trick Excel's UI uses) so they stay text, matching
# the source sheet where every Price/Volume cell is stored as a string.
$ws.Range("D2").Value = "26.305.27"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.689.94"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'219.03"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.5258"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "'0.2703"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "'0.06442"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'22.06"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'0.07473"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "1.689.50"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "'0.5856"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "'64.60"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "26.358.87"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "'4.957"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'10.89"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'189.80"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "'6.221"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'1.007"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "'144.86"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'7.667"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("E26").Value = "  +4.63%  "
$ws.Range("D27").Value = "'15.87"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "'0.06683"
$ws.Range("E28").Value = "  +14.65%  "
$ws.Range("D29").Value = "'1.351"
$ws.Range("E29").Value = "  +5.14%  "
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").Value = "'3.587"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").Value = "'0.6226"
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("D36").Value = "'2.392"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").Value = "'2.697"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("D38").Value = "'6.352"
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "1.107.54"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "'0.8868"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "'1.017"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "'100.99"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "1.837.08"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "'56.92"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.179"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "'0.05267"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "'0.4302"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "'6.060"
$ws.Range("E51").Value = "  +2.71%  "
